$wb = $excel.ActiveWorkbook

$data = @(
    @("Bitcoin", "BTC", 93593.70593424888, 1853167919054.723, 59748244183.16138, -2.16511306),
    @("Ethereum", "ETH", 3352.555608610053, 403839631695.9297, 32513038434.96328, 0.6943679699999999),
    @("Tether USDt", "USDT", 0.9986123728111663, 139654212669.495, 134856672030.5432, -0.08507927),
    @("XRP", "XRP", 2.190784119817687, 125427605845.8971, 10926323035.3287, -2.13208952),
    @("BNB", "BNB", 681.414910328032, 98128919741.89166, 1680950145.051784, 3.53900119),
    @("Solana", "SOL", 185.0009753963551, 88700600404.16412, 4520706652.315876, 1.24663073),
    @("Dogecoin", "DOGE", 0.3133869202417922, 46168767927.3466, 4122155365.7529, -1.0920213),
    @("USDC", "USDC", 0.9999344736902294, 43179223760.28349, 9172186703.4751, -0.01844355),
    @("Cardano", "ADA", 0.8959722778375632, 31472074199.14585, 1222973140.17203, 0.04080606),
    @("TRON", "TRX", 0.2498807197361438, 21543197412.80721, 1235717224.513316, 1.14241428),
    @("Avalanche", "AVAX", 37.15734698540833, 15228405331.77435, 637931089.3408355, -0.53231795),
    @("Chainlink", "LINK", 22.89203701503313, 14607408142.89674, 1084790623.173182, 3.00272092),
    @("Toncoin", "TON", 5.480738466563957, 13986289946.69558, 230554818.6848345, 0.92166012),
    @("Shiba Inu", "SHIB", 0.00002183443620001012, 12866061001.25375, 658972883.8270316, -0.20632331),
    @("Sui", "SUI", 4.269991333044081, 12501082905.34613, 1995715952.473691, -1.32131459),
    @("Polkadot", "DOT", 7.060517600388184, 10814144493.51163, 421613573.151455, 1.13592866),
    @("Stellar", "XLM", 0.3572769795416412, 10799733169.8659, 439595598.3194231, -1.06022475),
    @("Hedera", "HBAR", 0.2707210739885201, 10355666527.91892, 995252051.6471295, 1.31163355),
    @("Bitcoin Cash", "BCH", 441.961743447117, 8753385181.408186, 401518573.6439794, -1.69652879),
    @("UNUS SED LEO", "LEO", 9.421996910912126, 8710119987.361288, 1120666.48264045, 0.94497427),
    @("Uniswap", "UNI", 13.81215888022934, 8293967618.970984, 596263820.7778939, -1.68286208),
    @("Litecoin", "LTC", 102.1565761590558, 7696318059.579557, 740881840.7011844, 0.41680271),
    @("Pepe", "PEPE", 0.00001803631128995296, 7587693992.939088, 2133190220.55208, -0.33329076),
    @("Hyperliquid", "HYPE", 27.85843433477825, 7547088524.498374, 740055505.3636495, -9.852136570000001),
    @("NEAR Protocol", "NEAR", 5.176209265976165, 6304137124.600403, 931955322.5777456, 2.13008712),
    @("Ethena USDe", "USDe", 0.998519221181791, 5914790057.995737, 104646477.902974, -0.08632832),
    @("Bitget Token", "BGB", 4.082166220556925, 5715032708.779695, 118860920.2149833, -1.75967283),
    @("Dai", "DAI", 0.9996664115784819, 5363592873.11825, 161988535.0059164, -0.01436735),
    @("Aave", "AAVE", 350.1955302810284, 5259326840.086153, 1220682102.293655, 13.72400358),
    @("Aptos", "APT", 9.300876569441261, 5174906641.566488, 508002701.5491415, -0.76327685),
    @("Internet Computer", "ICP", 10.20095846321917, 4870933423.085492, 195140183.5938376, 1.75047501),
    @("Cronos", "CRO", 0.1574520449695216, 4183746569.616963, 36328265.03279613, -0.24622047),
    @("POL (ex-MATIC)", "POL", 0.4802348569379254, 4032267067.336287, 241914602.8810993, 0.10246063),
    @("Ethereum Classic", "ETC", 26.52573272841433, 3980482721.925358, 256630611.4218283, 0.88244808),
    @("Mantle", "MNT", 1.178787309293449, 3968790277.59795, 146169415.3803568, 0.18555388),
    @("VeChain", "VET", 0.04627761778294641, 3747794781.725382, 102573910.5960994, 0.24786643),
    @("Render", "RENDER", 7.149030474303392, 3701028118.391366, 449416600.1176429, 0.4262359),
    @("Monero", "XMR", 191.6995528684213, 3536232590.808321, 77588184.84335074, -0.31656531),
    @("Bittensor", "TAO", 459.4568611964735, 3391221687.252054, 211836599.103518, -0.51182856),
    @("MANTRA", "OM", 3.564293625131811, 3385324133.628084, 49614970.62176432, -3.19365097),
    @("Arbitrum", "ARB", 0.7733696020917132, 3255972617.45372, 501526645.6670145, 2.40430834),
    @("Artificial Superintelligence Alliance", "FET", 1.28553368624621, 3136659296.322258, 314626950.2999723, 0.78394303),
    @("Ethena", "ENA", 1.037312582644361, 3047105711.517812, 681892040.1650341, -0.61048477),
    @("Filecoin", "FIL", 4.963628554085775, 3046068636.197654, 386894734.0834174, 0.81052723),
    @("Algorand", "ALGO", 0.3591466503228555, 2993404290.289616, 418512884.1953453, 2.63779096),
    @("Kaspa", "KAS", 0.117764762669642, 2992529400.661585, 171098889.8951296, -1.22420974),
    @("Fantom", "FTM", 0.9788627356913631, 2744373664.983165, 406763196.502709, -1.61386279),
    @("OKB", "OKB", 45.47082453481779, 2728249472.089067, 2329667.10878335, 1.47674835),
    @("Virtuals Protocol", "VIRTUAL", 2.654291183886701, 2654291183.886701, 238087780.0917314, 11.10587863),
    @("Cosmos", "ATOM", 6.57679106239419, 2571092578.851387, 220750890.7668709, 0.10150516),
)

$liveSheet = $wb.Worksheets.Item("Live Data")
$topSheet = $wb.Worksheets.Item("Top 5")

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $liveSheet.Cells.Item($row, 1).Value = $rec[0]
    $liveSheet.Cells.Item($row, 2).Value = $rec[1]
    $liveSheet.Cells.Item($row, 3).Value = $rec[2]
    $liveSheet.Cells.Item($row, 4).Value = $rec[3]
    $liveSheet.Cells.Item($row, 5).Value = $rec[4]
    $liveSheet.Cells.Item($row, 6).Value = $rec[5]
}

for ($i = 0; $i -lt 5; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $topSheet.Cells.Item($row, 1).Value = $rec[0]
    $topSheet.Cells.Item($row, 2).Value = $rec[1]
    $topSheet.Cells.Item($row, 3).Value = $rec[2]
    $topSheet.Cells.Item($row, 4).Value = $rec[3]
    $topSheet.Cells.Item($row, 5).Value = $rec[4]
    $topSheet.Cells.Item($row, 6).Value = $rec[5]
}
